$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 10's formatting into the new row 11 (preserves the amount
# column's cell style) before overwriting the values for the new entry.
$ws.Rows("10").Copy()
$ws.Rows("11").Insert()

# B11: date-like label - force text storage (avoid auto date conversion),
# then drop back to the default (unstyled) cell format.
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "2019.9.11"
$ws.Range("B11").ClearFormats()

# C11: amount - keep the style copied from row 10.
$ws.Range("C11").Value = 500

# D11: description - reuses the "老师转账" text, default formatting.
$ws.Range("D11").ClearFormats()
$ws.Range("D11").Value = "老师转账"

# E11: running balance.
$ws.Range("E11").Value = 663

$ws.Range("H10").Select()
